$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Price / Volume(1h) cells hold numeric-looking text ("310.80", "0.64%").
# Force each cell to Text format *before* assigning so Excel stores the
# literal string instead of silently coercing it to a number/percentage.
$textCells = @(
    "D2", "E2", "D3", "E3", "D4", "E4", "D5", "E5", "D6", "E6",
    "D7", "E7", "D8", "E8", "D9", "E9", "D10", "E10", "D11", "E11",
    "D12", "E12", "D13", "E13", "D14", "E14", "D15", "E15", "D16", "E16",
    "D17", "E17", "D18", "E18", "D19", "E19", "D20", "E20", "D21", "E21",
    "D22", "E22", "D23", "E23", "D24", "E24", "E25", "D26", "E26", "D27",
    "E27", "D39", "E39", "D40", "E40", "D41", "E41", "D42", "E42", "D43",
    "E43", "D44", "E44", "D45", "E45", "D46", "E46", "D47", "E47", "D48",
    "E48", "D49", "E49", "D50", "E50", "D51", "E51"
)
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Coin / Link columns (plain text, safe to assign directly)
$ws.Range("B9").Value = "MXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("B18").Value = "BTSEToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"

# Price / Volume(1h) columns (assigned as text)
$ws.Range("D2").Value = "310.80"
$ws.Range("E2").Value = "0.64%"
$ws.Range("D3").Value = "37.16"
$ws.Range("E3").Value = "-2.64%"
$ws.Range("D4").Value = "5.123"
$ws.Range("E4").Value = "1.17%"
$ws.Range("D5").Value = "0.07769"
$ws.Range("E5").Value = "-1.84%"
$ws.Range("D6").Value = "4.378"
$ws.Range("E6").Value = "-0.42%"
$ws.Range("D7").Value = "8.220"
$ws.Range("E7").Value = "-0.35%"
$ws.Range("D8").Value = "1.877"
$ws.Range("E8").Value = "-7.33%"
$ws.Range("D9").Value = "0.9191"
$ws.Range("E9").Value = "-1.26%"
$ws.Range("D10").Value = "0.1216"
$ws.Range("E10").Value = "-4.93%"
$ws.Range("D11").Value = "0.1897"
$ws.Range("E11").Value = "-0.50%"
$ws.Range("D12").Value = "0.09169"
$ws.Range("E12").Value = "4.99%"
$ws.Range("D13").Value = "0.03435"
$ws.Range("E13").Value = "-0.63%"
$ws.Range("D14").Value = "0.09683"
$ws.Range("E14").Value = "-0.70%"
$ws.Range("D15").Value = "0.001365"
$ws.Range("E15").Value = "-2.96%"
$ws.Range("D16").Value = "0.005970"
$ws.Range("E16").Value = "-5.20%"
$ws.Range("D17").Value = "3.555"
$ws.Range("E17").Value = "-0.20%"
$ws.Range("D18").Value = "3.035"
$ws.Range("E18").Value = "-0.37%"
$ws.Range("D19").Value = "0.3406"
$ws.Range("E19").Value = "-1.04%"
$ws.Range("D20").Value = "5.260"
$ws.Range("E20").Value = "4.75%"
$ws.Range("D21").Value = "0.1268"
$ws.Range("E21").Value = "-2.08%"
$ws.Range("D22").Value = "0.2591"
$ws.Range("E22").Value = "2.80%"
$ws.Range("D23").Value = "0.02105"
$ws.Range("E23").Value = "5,591.10%"
$ws.Range("D24").Value = "0.04368"
$ws.Range("E24").Value = "0.80%"
$ws.Range("E25").Value = "-1.98%"
$ws.Range("D26").Value = "0.004250"
$ws.Range("E26").Value = "-8.18%"
$ws.Range("D27").Value = "0.0001301"
$ws.Range("E27").Value = "-63.79%"
$ws.Range("D39").Value = "0.02094"
$ws.Range("E39").Value = "-7.02%"
$ws.Range("D40").Value = "0.05035"
$ws.Range("E40").Value = "-0.39%"
$ws.Range("D41").Value = "0.007665"
$ws.Range("E41").Value = "1.00%"
$ws.Range("D42").Value = "0.009780"
$ws.Range("E42").Value = "-1.59%"
$ws.Range("D43").Value = "0.1346"
$ws.Range("E43").Value = "-1.72%"
$ws.Range("D44").Value = "0.002061"
$ws.Range("E44").Value = "-1.77%"
$ws.Range("D45").Value = "0.009579"
$ws.Range("E45").Value = "8.27%"
$ws.Range("D46").Value = "0.00006714"
$ws.Range("E46").Value = "3.05%"
$ws.Range("D47").Value = "0.00000000750"
$ws.Range("E47").Value = "-0.34%"
$ws.Range("D48").Value = "0.001200"
$ws.Range("E48").Value = "-0.34%"
$ws.Range("D49").Value = "0.002935"
$ws.Range("E49").Value = "-2.47%"
$ws.Range("D50").Value = "0.00002101"
$ws.Range("E50").Value = "-0.34%"
$ws.Range("D51").Value = "0.0002001"
$ws.Range("E51").Value = "-0.34%"

# Restore default styling (no explicit style) on the text-formatted cells
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
